# Update the forecast-error table:
#  - rows 2-5 (quarters Q0-Q3) get new, corrected metric values
#  - rows 6-11 (quarters Q4-Q9) are newly appended with their own metric values
# Column layout: A=label, B=ME, C=MAE, D=MSE, E=RMSE, F=SE, G=N

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  Label = "Q0"; B = 0.1121240238528402; C = 0.7618801882577865; D = 3.665753451068926;  E = 1.914615745017503;  F = 1.932450120671928;  G = 46 },
    @{ Row = 3;  Label = "Q1"; B = 0.2384546409039757; C = 0.7169237981258901; D = 3.429598611108179;  E = 1.851917549759756;  F = 1.857253653004148;  G = 45 },
    @{ Row = 4;  Label = "Q2"; B = 0.1554256882402988; C = 0.746261807642007;  D = 3.542524738013336;  E = 1.882159594193153;  F = 1.897416690034527;  G = 44 },
    @{ Row = 5;  Label = "Q3"; B = 0.2282236696078792; C = 0.7489212578072509; D = 3.600547596828088;  E = 1.897510895048586;  F = 1.906029573040004;  G = 43 },
    @{ Row = 6;  Label = "Q4"; B = 0.2046491862038776; C = 0.7536088930238598; D = 3.651661730967382;  E = 1.910932162837651;  F = 1.922972664121156;  G = 42 },
    @{ Row = 7;  Label = "Q5"; B = 0.2392227887939334; C = 0.7909065342463171; D = 3.777102677640282;  E = 1.943476955778041;  F = 1.95265767950657;   G = 41 },
    @{ Row = 8;  Label = "Q6"; B = 0.2155957076965697; C = 0.7869099249525766; D = 3.856332724670826;  E = 1.963754751660916;  F = 1.976749783912624;  G = 40 },
    @{ Row = 9;  Label = "Q7"; B = 0.239545290171405;  C = 0.8122047914116975; D = 3.974152750272352;  E = 1.993527714949645;  F = 2.004954792540174;  G = 39 },
    @{ Row = 10; Label = "Q8"; B = 0.2295130447590046; C = 0.80090761193645;   D = 4.057843671499801;  E = 2.01440901296132;   F = 2.028155615889948;  G = 38 },
    @{ Row = 11; Label = "Q9"; B = 0.189254522310732;  C = 0.7743474122443053; D = 3.997164790338612;  E = 1.999291071939904;  F = 2.017767317389937;  G = 37 }
)

# Use the formatting already present on the A2 label cell (bold, centered,
# bordered) for the newly appended label cells in column A.
$ws.Range("A2").Copy() | Out-Null

foreach ($entry in $data) {
    $r = $entry.Row

    if ($r -gt 5) {
        $ws.Range("A$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    }
    $ws.Range("A$r").Value = $entry.Label

    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}

$excel.CutCopyMode = 0
